$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> { column letter -> new numeric value } derived from the
# author's edit: erroneous "filled" values in columns D (res_c4) and E
# (total) of rows 2-22 and 29-71 are reset to 0, and rows 23, 24, 25, 26,
# 27, 28 and 72 receive their real survey results across columns A-AR.

$changes = @{
    2 = @{ "D" = 0; "E" = 0 }
    3 = @{ "D" = 0; "E" = 0 }
    4 = @{ "D" = 0; "E" = 0 }
    5 = @{ "D" = 0; "E" = 0 }
    6 = @{ "D" = 0; "E" = 0 }
    7 = @{ "D" = 0; "E" = 0 }
    8 = @{ "D" = 0; "E" = 0 }
    9 = @{ "D" = 0; "E" = 0 }
    10 = @{ "D" = 0; "E" = 0 }
    11 = @{ "D" = 0; "E" = 0 }
    12 = @{ "D" = 0; "E" = 0 }
    13 = @{ "D" = 0; "E" = 0 }
    14 = @{ "D" = 0; "E" = 0 }
    15 = @{ "D" = 0; "E" = 0 }
    16 = @{ "D" = 0; "E" = 0 }
    17 = @{ "D" = 0; "E" = 0 }
    18 = @{ "D" = 0; "E" = 0 }
    19 = @{ "D" = 0; "E" = 0 }
    20 = @{ "D" = 0; "E" = 0 }
    21 = @{ "D" = 0; "E" = 0 }
    22 = @{ "D" = 0; "E" = 0 }
    23 = @{ "A" = 16.98; "B" = 23.019999999999996; "D" = 4.0999999999999996; "E" = 44.1; "I" = 2.15; "M" = 3.7; "O" = 3; "Q" = 3.9; "R" = 4.2300000000000004; "S" = 1.65; "T" = 0.97; "U" = 1.82; "V" = 0.8; "W" = 1.2; "X" = 2.4; "Y" = 2.7; "Z" = 0.73; "AA" = 1; "AB" = 5.13; "AC" = 2.2000000000000002; "AD" = 1.42; "AE" = 1; "AG" = 10.1; "AH" = 2.2000000000000002; "AI" = 2.2999999999999998; "AK" = 1.5; "AL" = 0.4; "AM" = 0.4; "AN" = 0.4; "AO" = 0.4; "AP" = 0.4; "AR" = 0.6 }
    24 = @{ "C" = 14.600000000000001; "D" = 0; "E" = 14.600000000000001 }
    25 = @{ "A" = 10.1; "B" = 14.11; "D" = 5.3999999999999995; "E" = 29.61; "I" = 2.41; "K" = 2.5499999999999998; "M" = 2.5499999999999998; "R" = 2.59; "S" = 1.25; "T" = 0.3; "U" = 0.75; "V" = 0.7; "W" = 0.7; "X" = 1.3; "Z" = 0.83; "AA" = 0.4; "AB" = 4.4000000000000004; "AC" = 1.4; "AD" = 1.28; "AE" = 0.35; "AF" = 0.45; "AG" = 11.1; "AH" = 1.5; "AK" = 1.6; "AL" = 0.1; "AM" = 0.3; "AN" = 0.3; "AO" = 0.3; "AP" = 0.3; "AQ" = 0.1; "AR" = 2.4 }
    26 = @{ "C" = 12.6; "D" = 0; "E" = 12.6 }
    27 = @{ "A" = 17.45; "B" = 14.91; "D" = 8.1000000000000014; "E" = 40.46; "I" = 2.8; "M" = 3.05; "O" = 4; "Q" = 4.7; "R" = 2.9; "S" = 1.25; "T" = 0.63; "U" = 2.52; "V" = 0.7; "W" = 0.85; "Z" = 0.98; "AA" = 0.88; "AB" = 2.99; "AC" = 1.5; "AD" = 2.16; "AE" = 0.45; "AG" = 14.02; "AH" = 1.1000000000000001; "AI" = 1.1000000000000001; "AK" = 4.9000000000000004 }
    28 = @{ "C" = 16.22; "D" = 0; "E" = 16.22 }
    29 = @{ "D" = 0; "E" = 0 }
    30 = @{ "D" = 0; "E" = 0 }
    31 = @{ "D" = 0; "E" = 0 }
    32 = @{ "D" = 0; "E" = 0 }
    33 = @{ "D" = 0; "E" = 0 }
    34 = @{ "D" = 0; "E" = 0 }
    35 = @{ "D" = 0; "E" = 0 }
    36 = @{ "D" = 0; "E" = 0 }
    37 = @{ "D" = 0; "E" = 0 }
    38 = @{ "D" = 0; "E" = 0 }
    39 = @{ "D" = 0; "E" = 0 }
    40 = @{ "D" = 0; "E" = 0 }
    41 = @{ "D" = 0; "E" = 0 }
    42 = @{ "D" = 0; "E" = 0 }
    43 = @{ "D" = 0; "E" = 0 }
    44 = @{ "D" = 0; "E" = 0 }
    45 = @{ "D" = 0; "E" = 0 }
    46 = @{ "D" = 0; "E" = 0 }
    47 = @{ "D" = 0; "E" = 0 }
    48 = @{ "D" = 0; "E" = 0 }
    49 = @{ "D" = 0; "E" = 0 }
    50 = @{ "D" = 0; "E" = 0 }
    51 = @{ "D" = 0; "E" = 0 }
    52 = @{ "D" = 0; "E" = 0 }
    53 = @{ "D" = 0; "E" = 0 }
    54 = @{ "D" = 0; "E" = 0 }
    55 = @{ "D" = 0; "E" = 0 }
    56 = @{ "D" = 0; "E" = 0 }
    57 = @{ "D" = 0; "E" = 0 }
    58 = @{ "D" = 0; "E" = 0 }
    59 = @{ "D" = 0; "E" = 0 }
    60 = @{ "D" = 0; "E" = 0 }
    61 = @{ "D" = 0; "E" = 0 }
    62 = @{ "D" = 0; "E" = 0 }
    63 = @{ "D" = 0; "E" = 0 }
    64 = @{ "D" = 0; "E" = 0 }
    65 = @{ "D" = 0; "E" = 0 }
    66 = @{ "D" = 0; "E" = 0 }
    67 = @{ "D" = 0; "E" = 0 }
    68 = @{ "D" = 0; "E" = 0 }
    69 = @{ "D" = 0; "E" = 0 }
    70 = @{ "D" = 0; "E" = 0 }
    71 = @{ "D" = 0; "E" = 0 }
    72 = @{ "A" = 22.85; "B" = 17.64; "D" = 6.3000000000000007; "E" = 46.790000000000006; "I" = 4; "K" = 4; "M" = 2.65; "O" = 4; "Q" = 4.7; "R" = 3.5; "S" = 0.75; "T" = 0.43; "U" = 1.57; "V" = 1.0900000000000001; "W" = 0.8; "X" = 0.8; "Y" = 1.1000000000000001; "Z" = 0.93; "AA" = 0.85; "AB" = 3.49; "AC" = 3.1; "AD" = 1.1000000000000001; "AE" = 0.53; "AF" = 1.1000000000000001; "AG" = 13.47; "AH" = 2.2000000000000002; "AI" = 2.2999999999999998; "AK" = 2.1 }
}

foreach ($rowKey in $changes.Keys) {
    $rowData = $changes[$rowKey]
    foreach ($colKey in $rowData.Keys) {
        $ws.Range("$colKey$rowKey").Value = $rowData[$colKey]
    }
}

# Restore the active cell selection on the sheet to match the saved view.
$ws.Range("I9").Select()
